# Update a handful of imputed values in the RandomForest result sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D7").Value = -7.1948
$ws.Range("B10").Value = 5.045699999999998
$ws.Range("B12").Value = 4.6487
$ws.Range("D15").Value = -7.8929
$ws.Range("B18").Value = 7.179399999999998
$ws.Range("D20").Value = -7.542899999999992
$ws.Range("D29").Value = -7.384699999999999
$ws.Range("D30").Value = -7.303100000000001
$ws.Range("D31").Value = -8.668399999999997
$ws.Range("B37").Value = 9.070399999999999
$ws.Range("D40").Value = -7.884599999999995
$ws.Range("B55").Value = 5.823299999999998
$ws.Range("B68").Value = 4.529399999999997
$ws.Range("D68").Value = -7.200799999999998
$ws.Range("D76").Value = -7.327599999999995
$ws.Range("B77").Value = 9.061300000000003
$ws.Range("B78").Value = 9.583200000000001
$ws.Range("D87").Value = -7.721499999999996
$ws.Range("D88").Value = -7.119099999999996
$ws.Range("D96").Value = -7.448199999999999
$ws.Range("D98").Value = -8.552299999999999
$ws.Range("D101").Value = -7.819899999999999
$ws.Range("D102").Value = -7.798299999999998
